# Q2 answer: 97.3% -> 80.5%  (split into two runs: "80.5" and "%")
$d = $word.ActiveDocument

$f1 = $d.Content
$f1.Find.Execute("97.3%", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s1 = $f1.Start
$e1 = $f1.End
$r1 = $d.Range($s1, $e1)
$xml1 = '<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>80.5</w:t></w:r><w:r><w:t>%</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1)

# Q3 answer: "There are no changes ???" -> "83.5%" (split into two runs: "83.5" and "%"),
# preserving the _GoBack bookmark that trails the text.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

$f2 = $d.Content
$f2.Find.Execute("There are no changes ???", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s2 = $f2.Start
$e2 = $f2.End
$r2 = $d.Range($s2, $e2)
$xml2 = '<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>83.5</w:t></w:r><w:r><w:t>%</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml2)
